# tokyo nova typo fix
# Fix the English title typo "Toky Nova" -> "Tokyo Nova" in the checklist.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C5").Value = "Hounds Afternoon: Tokyo Nova The 2nd Edition Replay"

# Move the active selection to the edited cell, matching the authored workbook.
$ws.Range("C5").Select()
